# Add a "Hungary" market tab, cloned from the existing "Slovakia" tab
# (same layout/styles/merged cells), with its own market name and
# Jira/NGC reference string, per commit "Added HungaryFC Test data".

$wb = $excel.ActiveWorkbook

$slovakia = $wb.Worksheets.Item("Slovakia")

# Mimic Excel's "Move or Copy... > Create a copy" placed right after Slovakia.
$slovakia.Copy($null, $slovakia)
$hungary = $wb.Worksheets.Item($slovakia.Index + 1)
$hungary.Name = "Hungary"

# Fill in the Hungary-specific market name and NGC/Jira reference list.
$hungary.Range("B2").Value = "Hungary Market"
$hungary.Range("B4").Value = "NGC-4308/T3592/T3595/T3617/T3620"

# Leave the source (Slovakia) sheet with everything selected, deselected as
# the active tab, matching the state Excel leaves a sheet in right after a
# "Create a copy" operation where the user had selected all cells.
$slovakia.Select()
$slovakia.Cells.Select()

# The new Hungary sheet becomes the active tab, with B4 (the value just
# entered) as the selected cell.
$hungary.Select()
$hungary.Range("B4").Select()
